# "added 3rd design iteration"
# Adds a new BOM line (row 19) for a "Primary side splice" part, and
# updates the sheet's view state (zoom + scroll position + selection)
# to match where the author was working when they saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New BOM row: Item / MPN / DPN / Quantity / Description ---------------
$ws.Range("A19").Value = "Primary side splice"

$ws.Range("B19").Value = 8383
$ws.Range("B19").HorizontalAlignment = -4131   # xlLeft

$ws.Range("C19").Value = "36-8383-ND"
$ws.Range("D19").Value = 2
$ws.Range("E19").Value = "Terminal Inline, Tap Connector IDC 14-16 to 18-20 AWG Blue"

# --- View state: zoom out a touch, scroll right one column, and leave ------
# --- the selection on the new quantity cell --------------------------------
$excel.ActiveWindow.Zoom = 125
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1

$ws.Range("D21").Select()
